# This workbook is a weekly price log for "Ajo" (garlic) at "Feria Lagunitas
# de Puerto Montt". Each week's refresh inserts the newest daily price
# record(s) near the top of the data block (pushing the older rows down),
# matching the commit message "Fruta / hortaliza, semanal" (weekly).
#
# Two new observations are inserted this week:
#   - 2023-03-21 (serial 45006) inserted before the existing row 311
#   - 2023-03-20 (serial 45005) inserted before the (now shifted) old row 362
#
# Every other row keeps its original data, just shifted down by the
# corresponding number of inserted rows above it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insertion 1: new row at 311 ---------------------------------------
$ws.Rows.Item(311).Insert()

$ws.Cells.Item(311, 1).Value = 4
$ws.Cells.Item(311, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(311, 3).Value = "Los Lagos"
$ws.Cells.Item(311, 4).Value = 45006
$ws.Cells.Item(311, 5).Value = 10
$ws.Cells.Item(311, 6).Value = 100112003
$ws.Cells.Item(311, 7).Value = "Ajo"
$ws.Cells.Item(311, 8).Value = "Chino"
$ws.Cells.Item(311, 9).Value = "Primera"
$ws.Cells.Item(311, 10).Value = 240
$ws.Cells.Item(311, 11).Value = 20000
$ws.Cells.Item(311, 12).Value = 21000
$ws.Cells.Item(311, 13).Value = 20500
$ws.Cells.Item(311, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(311, 15).Value = "China"
$ws.Cells.Item(311, 16).Value = 2050
$ws.Cells.Item(311, 17).Value = 10
$ws.Cells.Item(311, 18).Value = "Hortaliza"

# --- Insertion 2: new row at 362 (after insertion 1 has already shifted
#     everything from the old row 311 onward down by one) ---------------
$ws.Rows.Item(362).Insert()

$ws.Cells.Item(362, 1).Value = 4
$ws.Cells.Item(362, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(362, 3).Value = "Los Lagos"
$ws.Cells.Item(362, 4).Value = 45005
$ws.Cells.Item(362, 5).Value = 10
$ws.Cells.Item(362, 6).Value = 100112003
$ws.Cells.Item(362, 7).Value = "Ajo"
$ws.Cells.Item(362, 8).Value = "Chino"
$ws.Cells.Item(362, 9).Value = "Primera"
$ws.Cells.Item(362, 10).Value = 80
$ws.Cells.Item(362, 11).Value = 20000
$ws.Cells.Item(362, 12).Value = 20000
$ws.Cells.Item(362, 13).Value = 20000
$ws.Cells.Item(362, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(362, 15).Value = "China"
$ws.Cells.Item(362, 16).Value = 2000
$ws.Cells.Item(362, 17).Value = 10
$ws.Cells.Item(362, 18).Value = "Hortaliza"
